$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every cell (even the numeric-looking "23") as
# text, so first mark the two new rows as Text. This makes the upcoming
# .Value assignments keep "23" as a string instead of being coerced into
# the number 23.
$ws.Range("A4:H5").NumberFormat = "@"

$ws.Range("A4").Value = ""
$ws.Range("B4").Value = 'حسن '
$ws.Range("C4").Value = "23"
$ws.Range("D4").Value = 'الصمود'
$ws.Range("E4").Value = 'الرحلة 2'
$ws.Range("F4").Value = "C1"
$ws.Range("G4").Value = "UNICEF"
$ws.Range("H4").Value = '٠١‏/٠٥‏/٢٠٢٥ ٠٦:٣٢:٤١ م'

$ws.Range("A5").Value = ""
$ws.Range("B5").Value = 'حسن '
$ws.Range("C5").Value = "23"
$ws.Range("D5").Value = 'الصمود'
$ws.Range("E5").Value = 'الرحلة 2'
$ws.Range("F5").Value = "C1"
$ws.Range("G5").Value = "UNICEF"
$ws.Range("H5").Value = '٠١‏/٠٥‏/٢٠٢٥ ٠٦:٣٢:٤١ م'

# The values are already committed as text at this point; drop the
# temporary Text number-format again so the new rows end up using the
# workbook's default (Normal) cell style, same as the existing rows.
$ws.Range("A4:H5").Style = "Normal"
